# This script reproduces the commit "textes : petit nettoyage de fichiers
# annexes (trad Elders 3)" which appends three new paragraphs right after
# the paragraph ending in "The four uncles, too, were remarkable!" (and
# right before the two pre-existing trailing blank paragraphs).

$d = $word.ActiveDocument

# Locate the paragraph that ends with the known anchor text.
$anchorText = "The four uncles, too, were remarkable!"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains($anchorText)) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Anchor paragraph not found"
}

# Create three new, empty paragraphs right after the anchor paragraph.
$rng = $d.Paragraphs.Item($targetIndex).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null

$rng = $d.Paragraphs.Item($targetIndex + 1).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null

$rng = $d.Paragraphs.Item($targetIndex + 2).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null

# Fill in each of the three new paragraphs with its own OOXML, matching
# the exact run/formatting layout of the target revision.
$xmlPara1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="P68B1DB1-Normal1"/><w:spacing w:before="0" w:after="113"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Gentium Book Plus" w:hAnsi="Gentium Book Plus"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Now there's not a single man left from that generation: all my uncles are dead. The men of my parents' generation are all dead, there are none left. There's only one of my aunts left. There are only people from my generation left, and only one of our aunts. </w:t></w:r></w:p>
'@

$xmlPara2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="P68B1DB1-Normal1"/><w:spacing w:before="0" w:after="113"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Gentium Book Plus" w:hAnsi="Gentium Book Plus"/><w:lang w:val="en-US"/></w:rPr><w:t>In the past, my grandmother, well, people in the village called her ‘M</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Gentium Book Plus" w:hAnsi="Gentium Book Plus"/><w:lang w:val="en-US"/></w:rPr><w:t>other</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Gentium Book Plus" w:hAnsi="Gentium Book Plus"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Butee /pv̩˩ʈʰɯ˧/’. They said: ‘Mother Butee, you are remarkable! Your children, how obedient they are! Your children will be on a level with our times, they will find a place in today's society! How remarkable! Our children, what misery! They're not obedient!” That's what they said to my grandmother. My grandmother was a remarkable person! That's what you can say about her! We're the older generation now! The old ones, they're gone. There's only one aunt left!</w:t></w:r></w:p>
'@

$xmlPara3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="P68B1DB1-Normal1"/><w:spacing w:before="0" w:after="113"/><w:rPr><w:rFonts w:ascii="Gentium Book Plus" w:hAnsi="Gentium Book Plus"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr></w:r></w:p>
'@

$d.Paragraphs.Item($targetIndex + 1).Range.InsertXML($xmlPara1)
$d.Paragraphs.Item($targetIndex + 2).Range.InsertXML($xmlPara2)
$d.Paragraphs.Item($targetIndex + 3).Range.InsertXML($xmlPara3)
